$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.240.34'
$ws.Range('E2').Value = '  +3.67%  '

# Row 3
$ws.Range('D3').Value = '1.814.78'
$ws.Range('E3').Value = '  +4.31%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -1.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.31%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.69%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4439'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.77%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3736'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.68%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.81'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.17%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07703'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.15%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.126'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.12%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.07%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.09%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.297'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.87%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.572'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.64%  '

# Row 16
$ws.Range('D16').Value = '1.820.47'
$ws.Range('E16').Value = '  +3.98%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +10.32%  '

# Row 18
$ws.Range('E18').Value = '  +3.06%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06513'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.80%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9997'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.79%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.93%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.252'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.72%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5340'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.76%  '

# Row 24
$ws.Range('D24').Value = '28.314.34'
$ws.Range('E24').Value = '  +3.69%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.20%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.050'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -14.98%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.07%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.29%  '

# Row 29
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').Value = '2.023.65'
$ws.Range('E29').Value = '  +3.80%  '

# Row 30
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.324'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.55%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.49'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.80%  '

# Row 32
$ws.Range('E32').Value = '  -6.21%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.862'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.47%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09223'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.41%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.676'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.74%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '13.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.27%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02349'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.91%  '

# Row 38
$ws.Range('E38').Value = '  +1.72%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.182'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.49%  '

# Row 40
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06213'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.70%  '

# Row 41
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6563'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.61%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.202'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.22%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.088'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.40%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9992'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.83%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.92%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.388'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.28%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6076'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.11%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.764'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.57%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '126.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.57%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.035'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.77%  '

# Row 51
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.151'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.26%  '
